$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 1
    3  = 0
    4  = 1
    5  = 1
    6  = 2
    7  = 2
    8  = 1
    9  = 0
    10 = 0
    11 = 0
    12 = 0
    13 = 1
    14 = 0
    15 = 0
    16 = 0
    17 = 0
    18 = 1
    19 = 1
    20 = 2
    21 = 1
    22 = 1
    23 = 1
    24 = 3
    25 = 1
    26 = 3
    27 = 1
    28 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
